# Horarios actualizados Linea 141 - 416
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (sheet1): new scrape at 03:23:38, totals 10 -> 13 rows,
# row 14 refreshed in place, two rows shift down (old row15 -> row16),
# and two brand-new arrivals appended as rows 17-18.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:23:38"
$ws1.Range("A3").Value = "Total filas: 13"

# Row 14: overwritten with the fresher scrape for 215A_EL PATO
$ws1.Range("A14").Value = "03:23:38"
$ws1.Range("B14").Value = "04:45"
$ws1.Range("D14").Value = 82

# Row 15: the previous row 14 content (215A_EL PATO / 110) re-appears here
$ws1.Range("A15").Value = "02:56:05"
$ws1.Range("B15").Value = "04:46"
$ws1.Range("C15").Value = "215A_EL PATO"
$ws1.Range("D15").Value = 110
$ws1.Range("E15").Value = "LP1912"

# Row 16: the old row 15 content, shifted down one row
$ws1.Range("A16").Value = "02:56:05"
$ws1.Range("B16").Value = "04:53"
$ws1.Range("C16").Value = "11_ETCHEVERRY"
$ws1.Range("D16").Value = 117
$ws1.Range("E16").Value = "LP1912"

# Row 17: brand-new arrival
$ws1.Range("A17").Value = "03:23:38"
$ws1.Range("B17").Value = "05:16"
$ws1.Range("C17").Value = "17_ROMERO"
$ws1.Range("D17").Value = 113
$ws1.Range("E17").Value = "LP1912"

# Row 18: brand-new arrival
$ws1.Range("A18").Value = "03:23:38"
$ws1.Range("B18").Value = "05:22"
$ws1.Range("C18").Value = "23_HERNANDEZ"
$ws1.Range("D18").Value = 119
$ws1.Range("E18").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (sheet2): same refresh, totals 4 -> 5 rows, one new
# row inserted before the old row 9 (which shifts to row 10).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:23:38"
$ws2.Range("A3").Value = "Total filas: 5"

# Row 9: brand-new, fresher scrape for 215A_EL PATO
$ws2.Range("A9").Value = "03:23:38"
$ws2.Range("B9").Value = "04:45"
$ws2.Range("C9").Value = "215A_EL PATO"
$ws2.Range("D9").Value = 82
$ws2.Range("E9").Value = "LP1912"

# Row 10: the old row 9 content, shifted down one row
$ws2.Range("A10").Value = "02:56:05"
$ws2.Range("B10").Value = "04:46"
$ws2.Range("C10").Value = "215A_EL PATO"
$ws2.Range("D10").Value = 110
$ws2.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (sheet3): only the "last updated" timestamp changes.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 03:23:38"
